# Fix Mitelmandb sites, NCI Thesaurus 25.09e
$wb = $excel.ActiveWorkbook

# The "compounds" sheet holds the per-source metadata table, including the
# NCI Thesaurus row (source_version in column E).
$ws = $wb.Worksheets.Item("compounds")
$ws.Activate()

# Row 3 corresponds to the "NCI Thesaurus" source (column A = "NCI Thesaurus",
# column F = "nci"); column E is "source_version". Bump it from 25.08d to 25.09e.
$ws.Range("E3").Value = "25.09e"

# Reflect the last-edited cell as the active selection, matching the saved view.
$ws.Range("E3").Select()
